$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) from row 2 through row 379: 45189 -> 45190
for ($r = 2; $r -le 379; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45189) {
        $cell.Value = 45190
    }
}

# Row 379 gains an explicit row height (customHeight flag) in the target file
$ws.Rows.Item(379).RowHeight = 15

# Add new row 380 data
$ws.Cells.Item(380, 1).Value = "A 44208-2023"
$ws.Cells.Item(380, 2).Value = 45188
$ws.Cells.Item(380, 3).Value = 45190
$ws.Cells.Item(380, 4).Value = "SÖDERMANLANDS LÄN"
$ws.Cells.Item(380, 5).Value = "GNESTA"
$ws.Cells.Item(380, 6).Value = "Holmen skog AB"
$ws.Cells.Item(380, 7).Value = 0.7
$ws.Cells.Item(380, 8).Value = 0
$ws.Cells.Item(380, 9).Value = 0
$ws.Cells.Item(380, 10).Value = 0
$ws.Cells.Item(380, 11).Value = 0
$ws.Cells.Item(380, 12).Value = 0
$ws.Cells.Item(380, 13).Value = 0
$ws.Cells.Item(380, 14).Value = 0
$ws.Cells.Item(380, 15).Value = 0
$ws.Cells.Item(380, 16).Value = 0
$ws.Cells.Item(380, 17).Value = 0

# Apply the same date number format used in column B/C (style s="1") to the new row
$ws.Cells.Item(380, 2).NumberFormat = $ws.Cells.Item(379, 2).NumberFormat
$ws.Cells.Item(380, 3).NumberFormat = $ws.Cells.Item(379, 3).NumberFormat

# R380 needs the same wrap-text style (s="2") used by the rest of column R
$ws.Cells.Item(379, 18).Copy()
$ws.Cells.Item(380, 18).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
